# Applies the metrics_20_7 update:
#  - Column A (model name) labels are reordered across rows 2..26
#  - Columns B..Q get one common new metrics row (same for every model row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of model names for rows 2..26 (row 1 is the header)
$newNames = @(
    "model_20_7_0",
    "model_20_7_22",
    "model_20_7_21",
    "model_20_7_20",
    "model_20_7_19",
    "model_20_7_18",
    "model_20_7_17",
    "model_20_7_16",
    "model_20_7_15",
    "model_20_7_14",
    "model_20_7_13",
    "model_20_7_23",
    "model_20_7_12",
    "model_20_7_10",
    "model_20_7_9",
    "model_20_7_8",
    "model_20_7_7",
    "model_20_7_6",
    "model_20_7_5",
    "model_20_7_4",
    "model_20_7_3",
    "model_20_7_2",
    "model_20_7_1",
    "model_20_7_11",
    "model_20_7_24"
)

# New common metric values for columns B..Q (identical on every data row)
# Values are parsed from strings because the scripting engine does not
# accept scientific-notation numeric literals directly.
$newMetrics = @(
    [double]"0.9999805300082542",
    [double]"0.9991182316315311",
    [double]"0.9999999999999197",
    [double]"0.9999693906937425",
    [double]"0.9999910533168559",
    [double]"1.817438266093276e-05",
    [double]"0.0008230920668076507",
    [double]"6.471593238127035e-14",
    [double]"2.380135856745718e-05",
    [double]"1.190067931608656e-05",
    [double]"0.0002731515751168351",
    [double]"0.004263142345844525",
    [double]"1.000035944600146",
    [double]"0.004444633141428372",
    [double]"95.83099500301532",
    [double]"140.9294005231387"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newNames[$i]
    for ($j = 0; $j -lt $newMetrics.Count; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $newMetrics[$j]
    }
}
